$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A152:D160").EntireRow.Copy($ws.Range("A162").EntireRow)
